$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1 gets a top+bottom border (-> new style borderId=4)
$ws1.Range("C1").Borders.Item(8).LineStyle = 1
$ws1.Range("C1").Borders.Item(9).LineStyle = 1

# D1 gets a top+bottom+right border (-> new style borderId=5)
$ws1.Range("D1").Borders.Item(8).LineStyle = 1
$ws1.Range("D1").Borders.Item(9).LineStyle = 1
$ws1.Range("D1").Borders.Item(10).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# C1 / F1 get a top+bottom border (-> new style borderId=4)
$ws2.Range("C1").Borders.Item(8).LineStyle = 1
$ws2.Range("C1").Borders.Item(9).LineStyle = 1
$ws2.Range("F1").Borders.Item(8).LineStyle = 1
$ws2.Range("F1").Borders.Item(9).LineStyle = 1

# D1 / G1 get a top+bottom+right border (-> new style borderId=5)
$ws2.Range("D1").Borders.Item(8).LineStyle = 1
$ws2.Range("D1").Borders.Item(9).LineStyle = 1
$ws2.Range("D1").Borders.Item(10).LineStyle = 1
$ws2.Range("G1").Borders.Item(8).LineStyle = 1
$ws2.Range("G1").Borders.Item(9).LineStyle = 1
$ws2.Range("G1").Borders.Item(10).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
